# Update cryptocurrency price/volume snapshot pulled by the scraper run.
# Rows 7-18 also shuffle the Coin/Link pair down by one rank (a new coin,
# GateToken, entered the top of that block), independent of the Price/
# Volume(1h) refresh applied to every row below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns hold plain text in this sheet (e.g.
# "335.78" / "44.10" / "1.77%"), so force text formatting before writing
# the new values -- otherwise COM auto-converts numeric-looking strings to
# real numbers/percentages and we lose the exact formatting (trailing
# zeros, "%" suffix, etc.).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2
Set-TextValue $ws.Range('D2') '335.26'
Set-TextValue $ws.Range('E2') '1.62%'

# Row 3
Set-TextValue $ws.Range('D3') '43.99'
Set-TextValue $ws.Range('E3') '6.21%'

# Row 4
Set-TextValue $ws.Range('D4') '5.759'
Set-TextValue $ws.Range('E4') '2.06%'

# Row 5
Set-TextValue $ws.Range('D5') '0.08383'
Set-TextValue $ws.Range('E5') '1.57%'

# Row 6
Set-TextValue $ws.Range('D6') '8.860'
Set-TextValue $ws.Range('E6') '1.15%'

# Row 7
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws.Range('D7') '4.522'
Set-TextValue $ws.Range('E7') '-0.30%'

# Row 8
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range('D8') '1.958'
Set-TextValue $ws.Range('E8') '-2.80%'

# Row 9
$ws.Range('B9').Value = 'BTSEToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue $ws.Range('D9') '2.879'
Set-TextValue $ws.Range('E9') '-3.08%'

# Row 10
$ws.Range('B10').Value = 'MXToken'
$ws.Range('C10').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D10') '0.9507'
Set-TextValue $ws.Range('E10') '2.81%'

# Row 11
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range('D11') '0.1247'
Set-TextValue $ws.Range('E11') '-2.06%'

# Row 12
$ws.Range('B12').Value = 'WazirX'
$ws.Range('C12').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range('D12') '0.1974'
Set-TextValue $ws.Range('E12') '0.80%'

# Row 13
$ws.Range('B13').Value = 'MandalaExchangeToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range('D13') '0.1017'
Set-TextValue $ws.Range('E13') '8.16%'

# Row 14
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range('D14') '0.04419'
Set-TextValue $ws.Range('E14') '12.53%'

# Row 15
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range('D15') '0.1067'
Set-TextValue $ws.Range('E15') '0.57%'

# Row 16
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range('D16') '0.001290'
Set-TextValue $ws.Range('E16') '-1.19%'

# Row 17
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws.Range('D17') '0.006049'
Set-TextValue $ws.Range('E17') '-1.10%'

# Row 18
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws.Range('D18') '3.495'
Set-TextValue $ws.Range('E18') '1.42%'

# Row 20
Set-TextValue $ws.Range('D20') '8.667'
Set-TextValue $ws.Range('E20') '3.80%'

# Row 21
Set-TextValue $ws.Range('D21') '0.1364'
Set-TextValue $ws.Range('E21') '-0.66%'

# Row 23
Set-TextValue $ws.Range('D23') '0.04415'
Set-TextValue $ws.Range('E23') '0.46%'

# Row 24
Set-TextValue $ws.Range('D24') '0.001256'
Set-TextValue $ws.Range('E24') '-0.10%'

# Row 25
Set-TextValue $ws.Range('D25') '0.004367'
Set-TextValue $ws.Range('E25') '1.14%'

# Row 26
Set-TextValue $ws.Range('D26') '0.0001262'
Set-TextValue $ws.Range('E26') '5.08%'

# Row 27
Set-TextValue $ws.Range('D27') '0.0003995'
Set-TextValue $ws.Range('E27') '-94.68%'

# Row 39
Set-TextValue $ws.Range('D39') '0.02838'
Set-TextValue $ws.Range('E39') '2.44%'

# Row 40
Set-TextValue $ws.Range('D40') '0.05938'
Set-TextValue $ws.Range('E40') '7.58%'

# Row 41
Set-TextValue $ws.Range('D41') '0.007938'
Set-TextValue $ws.Range('E41') '-0.05%'

# Row 42
Set-TextValue $ws.Range('D42') '0.1427'
Set-TextValue $ws.Range('E42') '0.45%'

# Row 43
Set-TextValue $ws.Range('D43') '0.008999'
Set-TextValue $ws.Range('E43') '0.66%'

# Row 44
Set-TextValue $ws.Range('D44') '0.002130'
Set-TextValue $ws.Range('E44') '-0.53%'

# Row 45
Set-TextValue $ws.Range('E45') '-13.97%'

# Row 46
Set-TextValue $ws.Range('D46') '0.00007228'
Set-TextValue $ws.Range('E46') '3.13%'

# Row 47
Set-TextValue $ws.Range('D47') '0.00000000751'
Set-TextValue $ws.Range('E47') '0.12%'

# Row 48
Set-TextValue $ws.Range('D48') '0.003202'
Set-TextValue $ws.Range('E48') '0.33%'

# Row 49
Set-TextValue $ws.Range('E49') '-0.30%'

# Row 50
Set-TextValue $ws.Range('D50') '0.00002103'
Set-TextValue $ws.Range('E50') '0.12%'

# Row 51
Set-TextValue $ws.Range('D51') '0.0002003'
Set-TextValue $ws.Range('E51') '0.12%'
